$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Rename "Margem de comercialização (total)*" -> "Margem de comercialização"
#    (rows 46:56, column B)
$ws.Range("B46:B56").Value = "Margem de comercialização"

# 2) Collapse the per-row "Pessoal ocupado em 31/12" .. "31/22" labels
#    (rows 57:67, column B) into a single repeated label
$ws.Range("B57:B67").Value = "Pessoal ocupado em 31/12"

# 3) Auto-fit column B so it is wide enough to show its longest label
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(2).ColumnWidth = 48.6

# 4) Update the view/selection state: select B57:B67 (active cell B57)
#    and scroll the window down so row 84 is visible near the top
$ws.Range("B57:B67").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 84
$excel.ActiveWindow.Panes.Item(1).ScrollRow = 84
